# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.082.04'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.659.51'
$ws.Range('E3').Value = '  -1.31%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.40'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5158'
$ws.Range('E6').Value = '  -2.61%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2575'
$ws.Range('E8').Value = '  -4.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06266'
$ws.Range('E9').Value = '  -0.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.87'
$ws.Range('E10').Value = '  -2.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07497'
$ws.Range('E11').Value = '  -0.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.657.69'
$ws.Range('E12').Value = '  -1.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.389'
$ws.Range('E13').Value = '  -2.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5371'
$ws.Range('E14').Value = '  -5.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '66.01'
$ws.Range('E15').Value = '  -0.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅7866'
$ws.Range('E16').Value = '  -3.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.079.45'
$ws.Range('E17').Value = '  -0.75%  '
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.675'
$ws.Range('E19').Value = '  -3.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '186.51'
$ws.Range('E20').Value = '  -1.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.14'
$ws.Range('E21').Value = '  -3.93%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.161'
$ws.Range('E22').Value = '  -1.05%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '148.07'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1207'
$ws.Range('E25').Value = '  -4.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.363'
$ws.Range('E26').Value = '  -3.58%  '
$ws.Range('E27').Value = '  -2.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.377'
$ws.Range('E28').Value = '  +1.98%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.06085'
$ws.Range('E29').Value = '  -5.29%  '
$ws.Range('E30').Value = '  -2.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.455'
$ws.Range('E31').Value = '  -2.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.392'
$ws.Range('E32').Value = '  -2.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.623'
$ws.Range('E33').Value = '  -2.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9817'
$ws.Range('E34').Value = '  -3.26%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.386'
$ws.Range('E35').Value = '  -1.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.749'
$ws.Range('E36').Value = '  +1.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5845'
$ws.Range('E37').Value = '  -4.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.101.89'
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01592'
$ws.Range('E39').Value = '  -1.57%  '
$ws.Range('E40').Value = '  -3.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8477'
$ws.Range('E41').Value = '  -2.53%  '
$ws.Range('E42').Value = '  -0.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.92'
$ws.Range('E43').Value = '  -0.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.811.11'
$ws.Range('E44').Value = '  -1.18%  '
$ws.Range('E45').Value = '  +0.33%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.000'
$ws.Range('E46').Value = '  -0.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.69'
$ws.Range('E47').Value = '  -4.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.997'
$ws.Range('E48').Value = '  -0.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05228'
$ws.Range('E49').Value = '  -0.65%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4242'
$ws.Range('E50').Value = '  -0.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.840'
$ws.Range('E51').Value = '  -2.35%  '
